# Gestion projet S2 - update row 17 (first IHMs) and refresh selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 17: fill in hours spent on "G/H/I" (réalisé) and "K/L/M" (restant)
# columns for the newly started task. Dependent formula cells
# (J17, N17, O17, W17, and the totals in row 24/25) recalc automatically.
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 2

$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 4

# Update the active selection on the sheet to reflect where the author
# left off editing.
$ws.Range("B19").Select()
